$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '69.658.76'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +0.34%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.499.13'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '576.00'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.04%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '166.79'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -1.14%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '2.497.84'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('E11').Value = '  +0.26%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.356'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +3.45%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '4.95'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +1.90%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '2.956.52'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.20%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '69.523.44'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('E17').Value = '  -0.25%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '2.504.38'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.01%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '11.19'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -1.06%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '7.46'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -4.47%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '348.37'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +0.11%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '3.91'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.94%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '1.94'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +0.92%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -0.08%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '70.32'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('E26').Value = '  +0.08%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '8.78'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -1.02%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.625.74'
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -0.90%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '7.80'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -0.76%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '459.31'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('E33').Value = '  -2.56%  '
$ws.Range('E34').Value = '  -0.82%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  +0.00%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '156.89'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('E43').Value = '  +0.12%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '38.17'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('E45').Value = '  -4.06%  '
$ws.Range('E46').Value = '  -5.50%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '141.16'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E49').Value = '  -1.59%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0732'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +0.15%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.578'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.56%  '

Write-Host "Applied 77 cell updates"
